$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers): swap the last two headers so H1=SamplePortion, I1=Result ---
$ws.Range("H1").Value = "SamplePortion"
$ws.Range("I1").Value = "Result"

# --- Row 2 (type row): the #float type now carries a unit annotation ---
$ws.Range("H2").Value = "#float,  unit:mg"
$ws.Range("I2").Value = "#float,  unit:mg"

# --- Row 3 (new): per-column description / mapping keys ---
$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"
$ws.Range("H3").Value = "#PriseEssai"
$ws.Range("I3").Value = "#Resultat"
